$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header "PCA - XGBoost Model" in D1, matching style of A1/B1/C1
$ws.Range("D1").Value = "PCA - XGBoost Model"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in Random Forest model stats for column C (rows 2-7)
$ws.Range("C2").Value = 0.4537
$ws.Range("C3").Value = 0.5604
$ws.Range("C4").Value = 0.7486
$ws.Range("C5").Value = 1.5109
$ws.Range("C6").Value = 0.2751
$ws.Range("C7").Value = 0.5158

# Copy number styling for C2:C6 from the matching B column cells
$ws.Range("B2").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("B5").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("B6").Copy()
$ws.Range("C6").PasteSpecial(-4122)

# C7 gets the built-in "Percent" cell style
$ws.Range("C7").Style = "Percent"
$ws.Range("C7").Value = 0.5158

# Adjust selection to D7, mirroring final cursor placement
$ws.Range("D7").Select()

# Resize columns to fit their new contents (best-fit-style autosize)
$ws.Columns("A").ColumnWidth = 45.666666666666664
$ws.Columns("B").ColumnWidth = 27.666666666666668
$ws.Columns("C").ColumnWidth = 34.833333333333336
$ws.Columns("D").ColumnWidth = 26.666666666666668

# Minor row-height re-layout (consistent with row re-rendering after edit)
$ws.Rows("2").RowHeight = 29.1
$ws.Rows("4").RowHeight = 27.95
$ws.Rows("5").RowHeight = 30.6
$ws.Rows("6").RowHeight = 33.95
$ws.Rows("7").RowHeight = 36.95
